$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (columns D, L, M, N, O, P, R, S), the rest of the row
# (A,B,C,E,F,G,H,I,J,K,Q,T) stays unchanged.
$rows = @{
    2 = @{ D = 44187; L = "Primera"; M = 350; N = 16000; O = 16000; P = 16000; R = "Región Metropolitana";   S = 1067 }
    3 = @{ D = 44187; L = "Segunda"; M = 300; N = 13000; O = 13000; P = 13000; R = "Región Metropolitana";   S = 867  }
    4 = @{ D = 44176; L = "Segunda"; M = 500; N = 15000; O = 16000; P = 15500; R = "Región Metropolitana";   S = 1033 }
    5 = @{ D = 44162; L = "Tercera"; M = 500; N = 15000; O = 16000; P = 15500; R = "Región de O'Higgins";    S = 1033 }
    6 = @{ D = 44194; L = "Segunda"; M = 300; N = 15000; O = 16000; P = 15500; R = "Región Metropolitana";   S = 1033 }
    7 = @{ D = 44159; L = "Tercera"; M = 400; N = 15500; O = 16000; P = 15750; R = "Región de O'Higgins";    S = 1050 }
    8 = @{ D = 44169; L = "Segunda"; M = 500; N = 15000; O = 16000; P = 15500; R = "Región de O'Higgins";    S = 1033 }
    9 = @{ D = 44166; L = "Segunda"; M = 600; N = 16000; O = 17000; P = 16500; R = "Región de O'Higgins";    S = 1100 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 12).Value = $vals.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value = $vals.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $vals.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 18).Value = $vals.R   # R: Origen
    $ws.Cells.Item($r, 19).Value = $vals.S   # S: Precio $/Kg
}
